# Weekly update: insert a new daily price record as the new top row (row 82)
# for "Choclo" / "Dulce o Americano" in Provincia de Limarí, shifting all the
# existing rows from 82:103 down to 83:104.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single blank row at row 82 - this shifts rows 82:103 down to 83:104
# and Excel carries the D-column (date) cell style down from the row above,
# exactly like the other rows in this block.
$ws.Rows("82:82").Insert()

# Populate the newly-inserted row 82 with the new record's data.
$ws.Range("A82").Value = 7
$ws.Range("B82").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C82").Value = "Ñuble"
$ws.Range("D82").Value = 44524
$ws.Range("E82").Value = 16
$ws.Range("F82").Value = 100112024
$ws.Range("G82").Value = "Choclo"
$ws.Range("H82").Value = "Dulce o Americano"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 60
$ws.Range("K82").Value = 16000
$ws.Range("L82").Value = 17000
$ws.Range("M82").Value = 16500
$ws.Range("N82").Value = "$/malla 60 unidades"
$ws.Range("O82").Value = "Provincia de Limarí"
$ws.Range("P82").Value = 275
$ws.Range("Q82").Value = 60
$ws.Range("R82").Value = "Hortaliza"
